$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the invalid end_date typo for the Code4Lib row.
$ws.Range("C14").Value = "2024-05-31"

# Rename / shorten task names in column A (rows 3-10, 12, 14-16).
# Row 2, 11, 13 are left untouched.
$ws.Range("A3").Value  = "student reserves - Phase I"
$ws.Range("A4").Value  = "student reserves - Phase II"
$ws.Range("A5").Value  = "student reserves - Phase II"
$ws.Range("A6").Value  = "Ruby on Rails Udemy"
$ws.Range("A7").Value  = "ASpace PUI streamlining"
$ws.Range("A8").Value  = "ELUNA presentation"
$ws.Range("A9").Value  = "LOCKSS upgrade and migration"
$ws.Range("A10").Value = "OASIS"
$ws.Range("A12").Value = "Leganto working group NERS"
$ws.Range("A14").Value = "loan rule Code4Lib article"
$ws.Range("A15").Value = "Read by QxMD"
$ws.Range("A16").Value = "SpineOMatic Alma cloud app"

# Move the active selection to A18 (matches the saved view state).
$ws.Range("A18").Select()
